# Weekly refresh of "Hortaliza, Terminal Hortofrutícola Agro Chillán - Betarraga" data.
# Two brand-new weekly observations (Primera / Segunda quality) are inserted at the
# top of the rolling data block (rows 434-435), every existing row in the block
# (old rows 434-535) shifts down by two positions (to rows 436-537), and the two
# rows that fall off the bottom of the original block (old rows 536-537) are
# appended as brand-new rows 538-539.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the existing data block (rows 434-537, all columns A-R) before overwriting
# anything - this captures old rows 434..535 (which become new rows 436..537) as well
# as old rows 536..537 (which become the newly appended rows 538..539).
$oldBlock = $ws.Range("A434:R537").Value2

# Shift the whole block down by two rows in one bulk write.
$ws.Range("A436:R539").Value2 = $oldBlock

# The two rows appended at the very bottom (538-539) are brand new cells, so they
# don't inherit the date style used throughout column D; restore it explicitly.
$ws.Range("D538:D539").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Write the two brand-new weekly entries into rows 434 and 435 (columns A, B, C, E,
# F, G, H, N, Q, R are constant across the whole data block and are left untouched).
$ws.Range("D434").Value2 = 44964
$ws.Range("I434").Value2 = "Primera"
$ws.Range("J434").Value2 = 300
$ws.Range("K434").Value2 = 700
$ws.Range("L434").Value2 = 800
$ws.Range("M434").Value2 = 750
$ws.Range("O434").Value2 = "Provincia de Diguillín"
$ws.Range("P434").Value2 = 150

$ws.Range("D435").Value2 = 44964
$ws.Range("I435").Value2 = "Segunda"
$ws.Range("J435").Value2 = 300
$ws.Range("K435").Value2 = 500
$ws.Range("L435").Value2 = 600
$ws.Range("M435").Value2 = 550
$ws.Range("O435").Value2 = "Provincia de Diguillín"
$ws.Range("P435").Value2 = 110
